$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.269.07"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "3.674.43"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'675.61"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").Value = "'157.92"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("E10").Value = "  -5.61%  "

$ws.Range("E11").Value = "  -2.49%  "

$ws.Range("E12").Value = "  -3.36%  "

$ws.Range("D13").Value = "4.293.25"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "'32.31"
$ws.Range("E14").Value = "  -3.96%  "

$ws.Range("D15").Value = "3.675.70"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "69.208.93"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").Value = "'16.05"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("D20").Value = "'467.35"
$ws.Range("E20").Value = "  -3.27%  "

$ws.Range("D21").Value = "'9.99"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").Value = "'0.649"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").Value = "'79.72"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").Value = "3.818.52"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -6.85%  "

$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  -5.41%  "

$ws.Range("D29").Value = "'2.68"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("E30").Value = "  -5.92%  "

$ws.Range("E31").Value = "  -3.75%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("E33").Value = "  -4.99%  "

$ws.Range("D34").Value = "'26.85"
$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("D35").Value = "3.667.07"
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("E36").Value = "  -4.76%  "

$ws.Range("D37").Value = "'8.15"
$ws.Range("E37").Value = "  -4.26%  "

$ws.Range("D38").Value = "'6.28"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").Value = "'174.71"
$ws.Range("E42").Value = "  +8.13%  "

$ws.Range("E43").Value = "  -4.40%  "

$ws.Range("D44").Value = "'0.939"

$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.70"
$ws.Range("E46").Value = "  -5.41%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'28.01"
$ws.Range("E47").Value = "  -7.43%  "

$ws.Range("D48").Value = "'0.000276"
$ws.Range("E48").Value = "  -4.42%  "

$ws.Range("E49").Value = "  -5.09%  "

$ws.Range("E50").Value = "  -3.85%  "

$ws.Range("D51").Value = "'7.78"
$ws.Range("E51").Value = "  -3.10%  "
